$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C4 and C6 change from text "-" (shared string) to numeric 0
$ws.Range("C4").Value = 0
$ws.Range("C6").Value = 0

# Move selection from Q5 to D8
$ws.Range("D8").Select()
